$d = $word.ActiveDocument

$targetText = "Sets Model Statements: (Transform, Mapping, Statement, Kind);"
$newListText = "Sets HATEOAS Model Statements: (Transform, Mapping, Statement, Kind);"
$replacementText = "Core Model: Assert Transform / Query Mapping of Statement Kind Object T. Resulting U : Flows Transforms Statements Kinds matching domains."

# Locate the existing "Sets Model Statements: ..." list paragraph.
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ($targetText + [char]13)) {
        $targetIndex = $i
        break
    }
}

if ($targetIndex -eq -1) {
    throw "Could not locate paragraph with text '$targetText'"
}

# Insert a brand-new list paragraph directly above it by appending a
# paragraph right after the *previous* paragraph; this way the new
# paragraph picks up the previous paragraph's (plain) formatting instead of
# inheriting the <w:u w:val="none"/> run-properties override that belongs
# to the "Sets Model Statements" paragraph itself.
$prevPara = $d.Paragraphs.Item($targetIndex - 1)
$prevPara.Range.InsertParagraphAfter()

# The blank paragraph just created now sits at $targetIndex, and the
# original "Sets Model Statements" paragraph has shifted to $targetIndex + 1.
$newPara = $d.Paragraphs.Item($targetIndex)
$newRange = $d.Range($newPara.Range.Start, $newPara.Range.End - 1)
$newRange.Text = $newListText

# Update the wording of the original paragraph (now one position later).
$existingPara = $d.Paragraphs.Item($targetIndex + 1)
$existingRange = $d.Range($existingPara.Range.Start, $existingPara.Range.End - 1)
$existingRange.Text = $replacementText
